# loads_1pm.xlsx: rename the header cell A1 from "Loads" to "Load Name"
# (per commit message: "Change of name of _about_ to README" / general
# doc & naming clean-up — here the column header "Loads" becomes
# "Load Name"). All other cell contents are left untouched; Excel's
# shared-string table is recompacted automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Load Name"

# Leave the selection parked on A1 (closest reachable approximation of
# the saved file's view state).
$ws.Range("A1").Select()
